$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-05 Sunday" "2025-01-06 Monday"

Replace-Text "535÷3=" "517÷7="
Replace-Text "364÷2=" "844÷7="
Replace-Text "564÷9=" "860÷5="
Replace-Text "731÷3=" "615÷5="
Replace-Text "289÷2=" "927÷8="
Replace-Text "936÷5=" "649÷2="
Replace-Text "973÷3=" "153÷2="
Replace-Text "485÷2=" "331÷7="
Replace-Text "840÷9=" "446÷2="
Replace-Text "427÷8=" "464÷8="
Replace-Text "282÷6=" "951÷3="
Replace-Text "135÷3=" "745÷8="
Replace-Text "190÷8=" "207÷2="
Replace-Text "852÷8=" "445÷5="
Replace-Text "440÷3=" "126÷7="
Replace-Text "813÷9=" "654÷8="
Replace-Text "425÷4=" "762÷6="
Replace-Text "803÷9=" "393÷2="
Replace-Text "842÷2=" "394÷6="
Replace-Text "586÷7=" "510÷7="
Replace-Text "647÷5=" "688÷2="
Replace-Text "482÷8=" "782÷4="
Replace-Text "747÷5=" "585÷8="
Replace-Text "648÷7=" "199÷9="
Replace-Text "310÷5=" "407÷7="
